$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the old merged ranges so the rebuild starts clean ---
$ws.Range("A4:L4").UnMerge()
$ws.Range("A7:L7").UnMerge()

# --- Clear all existing cell content/formatting ---
$ws.Cells.Clear()

# --- Write the new cell values (this also extends the sheet's dimension) ---
    $ws.Range("A1").Value = "MASTER PACKAGE"
    $ws.Range("B1").Value = ""
    $ws.Range("C1").Value = ""
    $ws.Range("D1").Value = ""
    $ws.Range("A2").Value = "WesternGlove Centric8 PROD"
    $ws.Range("B2").Value = "M12225BVS563:KONRAD"
    $ws.Range("C2").Value = "BOM"
    $ws.Range("D2").Value = "MASTER"
    $ws.Range("A3").Value = "Placements"
    $ws.Range("B3").Value = ""
    $ws.Range("C3").Value = ""
    $ws.Range("D3").Value = ""
    $ws.Range("A5").Value = "Code"
    $ws.Range("B5").Value = "Product"
    $ws.Range("C5").Value = "Type"
    $ws.Range("D5").Value = "Description :"
    $ws.Range("E5").Value = "Main" + [char]10 + "Material"
    $ws.Range("F5").Value = "Composition"
    $ws.Range("G5").Value = "Coating" + [char]10 + "Composition"
    $ws.Range("H5").Value = "DUNE WASH"
    $ws.Range("I5").Value = "Weight" + [char]10 + "/ Yield"
    $ws.Range("J5").Value = "Common" + [char]10 + "Qty"
    $ws.Range("K5").Value = "Image"
    $ws.Range("L5").Value = "Supplier"
    $ws.Range("A6").Value = "SJC-" + [char]10 + "32234HOOAS"
    $ws.Range("B6").Value = "BLACK 3`"" + [char]10 + "SWIFT TACK"
    $ws.Range("C6").Value = "Accessory"
    $ws.Range("D6").Value = ""
    $ws.Range("E6").Value = ""
    $ws.Range("F6").Value = ""
    $ws.Range("G6").Value = ""
    $ws.Range("H6").Value = ""
    $ws.Range("I6").Value = ""
    $ws.Range("J6").Value = "1"
    $ws.Range("K6").Value = ""
    $ws.Range("L6").Value = "REVOLUTION" + [char]10 + "GROUP (HK)"
    $ws.Range("A7").Value = ""
    $ws.Range("B7").Value = "CLEAR LEG" + [char]10 + "STICKER"
    $ws.Range("C7").Value = "Accessory"
    $ws.Range("D7").Value = "SET ON RIGHT" + [char]10 + "BACK LEG PANEL" + [char]10 + "(AS WORN) 3/4`" IN" + [char]10 + "FROM FINISHED" + [char]10 + "SIDE SEAM. *SET" + [char]10 + "AFTER WASH*"
    $ws.Range("E7").Value = ""
    $ws.Range("F7").Value = ""
    $ws.Range("G7").Value = ""
    $ws.Range("H7").Value = "MENS - SJC-"
    $ws.Range("I7").Value = ""
    $ws.Range("J7").Value = ""
    $ws.Range("K7").Value = ""
    $ws.Range("L7").Value = "R-Pac" + [char]10 + "International" + [char]10 + "Corp"
    $ws.Range("A8").Value = "A La Carte Menu Up-Charges (For internal use - predicting BOM cost roll-up) (2)"
    $ws.Range("B8").Value = ""
    $ws.Range("C8").Value = ""
    $ws.Range("D8").Value = ""
    $ws.Range("E8").Value = ""
    $ws.Range("F8").Value = ""
    $ws.Range("G8").Value = ""
    $ws.Range("H8").Value = ""
    $ws.Range("I8").Value = ""
    $ws.Range("J8").Value = ""
    $ws.Range("K8").Value = ""
    $ws.Range("L8").Value = ""
    $ws.Range("A9").Value = ""
    $ws.Range("B9").Value = "base" + [char]10 + "cost"
    $ws.Range("C9").Value = "Special" + [char]10 + "CM"
    $ws.Range("D9").Value = ""
    $ws.Range("E9").Value = ""
    $ws.Range("F9").Value = ""
    $ws.Range("G9").Value = ""
    $ws.Range("H9").Value = ""
    $ws.Range("I9").Value = ""
    $ws.Range("J9").Value = ""
    $ws.Range("K9").Value = ""
    $ws.Range("L9").Value = ""
    $ws.Range("A10").Value = ""
    $ws.Range("B10").Value = "MENS OTHER" + [char]10 + "TRIM" + [char]10 + "DETAILS"
    $ws.Range("C10").Value = "Special" + [char]10 + "CM"
    $ws.Range("D10").Value = ""
    $ws.Range("E10").Value = ""
    $ws.Range("F10").Value = ""
    $ws.Range("G10").Value = ""
    $ws.Range("H10").Value = ""
    $ws.Range("I10").Value = ""
    $ws.Range("J10").Value = ""
    $ws.Range("K10").Value = ""
    $ws.Range("L10").Value = "SILVERMOON" + [char]10 + "JEANS" + [char]10 + "LIMITED"
    $ws.Range("A11").Value = "Displaying 17 - 20 of 20 results"
    $ws.Range("B11").Value = ""
    $ws.Range("C11").Value = ""
    $ws.Range("D11").Value = ""
    $ws.Range("E11").Value = ""
    $ws.Range("F11").Value = ""
    $ws.Range("G11").Value = ""
    $ws.Range("H11").Value = ""
    $ws.Range("I11").Value = ""
    $ws.Range("J11").Value = ""
    $ws.Range("K11").Value = ""
    $ws.Range("L11").Value = ""


# --- Re-apply the bordered / wrapped / top-left-aligned cell style used
#     throughout this sheet (matches cell style index s="1") to every cell
#     we just populated, A1:L11 ---
$fmtRange = $ws.Range("A1:L11")
$fmtRange.Borders.LineStyle = 1
$fmtRange.WrapText = $true
$fmtRange.HorizontalAlignment = -4131
$fmtRange.VerticalAlignment = -4160

# --- Column B needs to be widened (bumped by 10 characters) ---
$ws.Range("B1").ColumnWidth = 21.6

# --- Re-create the merged banner rows at their new positions ---
$ws.Range("A8:L8").Merge()
$ws.Range("A11:L11").Merge()

# --- Rename the sheet ---
$ws.Name = "Sheet1"
